# Daily attendance processing - reorder "Recorded By" (column G) names so that
# "System" (or, if "System" absent, "admin@admin.com") is listed first among
# the comma-separated names, preserving the relative order of the rest.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -eq $null -or $val -eq "") {
        continue
    }

    $parts = $val -split ',\s*'
    if ($parts.Count -lt 2) {
        continue
    }

    $hasSystem = $false
    $hasAdmin = $false
    foreach ($p in $parts) {
        if ($p.Equals("System")) { $hasSystem = $true }
        if ($p.Equals("admin@admin.com")) { $hasAdmin = $true }
    }

    if ($hasSystem) {
        $lead = "System"
    } elseif ($hasAdmin) {
        $lead = "admin@admin.com"
    } else {
        continue
    }

    $rest = @()
    $removed = $false
    foreach ($p in $parts) {
        if ((-not $removed) -and $p.Equals($lead)) {
            $removed = $true
        } else {
            $rest += $p
        }
    }

    $newVal = $lead + ", " + ($rest -join ", ")
    $cell.Value = $newVal
}
